$d = $word.ActiveDocument
$found = $d.Content.Find.Execute("Professional Experience", $false, $false, $false, $false, $false, $true, 1, $false, "Professional Experience", 2)
Write-Output "found=$found"
